$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1426.6666
$ws.Range("I62").Value = 1426.6666
$ws.Range("K62").Value = 1426.6666
$ws.Range("M62").Value = -802.6666

# Row 65
$ws.Range("H65").Value = 1426.6666
$ws.Range("I65").Value = 1426.6666
$ws.Range("K65").Value = 7133.333000000001
$ws.Range("M65").Value = -4013.333000000001

# Row 86
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# Row 98
$ws.Range("H98").Value = 36104.31
$ws.Range("I98").Value = 51190.25
$ws.Range("J98").Value = 2580
$ws.Range("K98").Value = 51190.25
$ws.Range("L98").Value = 2580
$ws.Range("M98").Value = -49692.25
$ws.Range("N98").Value = -5576

# Row 122
$ws.Range("H122").Value = 36104.31
$ws.Range("I122").Value = 51190.25
$ws.Range("J122").Value = 2580
$ws.Range("K122").Value = 153570.75
$ws.Range("L122").Value = 7740
$ws.Range("M122").Value = -151120.75
$ws.Range("N122").Value = -12640

# Row 132
$ws.Range("H132").Value = 2934304.5
$ws.Range("I132").Value = 3368276.5
$ws.Range("K132").Value = 10104829.5
$ws.Range("M132").Value = -10102299.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1323.0741
$ws.Range("I61").Value = 1116.421
$ws.Range("J61").Value = 1813.875
$ws.Range("K61").Value = 1116.421
$ws.Range("L61").Value = 1813.875
$ws.Range("M61").Value = -904.421
$ws.Range("N61").Value = -2237.875

# Row 74
$ws.Range("H74").Value = 26316870
$ws.Range("I74").Value = 41667500
$ws.Range("J74").Value = 1502
$ws.Range("K74").Value = 41667500
$ws.Range("L74").Value = 1502
$ws.Range("M74").Value = -41666626
$ws.Range("N74").Value = -3250

# Row 77
$ws.Range("H77").Value = 26316870
$ws.Range("I77").Value = 41667500
$ws.Range("J77").Value = 1502
$ws.Range("K77").Value = 208337500
$ws.Range("L77").Value = 7510
$ws.Range("M77").Value = -208333132
$ws.Range("N77").Value = -16246

# Row 132
$ws.Range("H132").Value = 1525.9474
$ws.Range("I132").Value = 1177.8572
$ws.Range("J132").Value = 1955.9412
$ws.Range("K132").Value = 3533.5716
$ws.Range("L132").Value = 5867.8236
$ws.Range("M132").Value = -1003.5716
$ws.Range("N132").Value = -10927.8236

# Row 135
$ws.Range("H135").Value = 23285.8
$ws.Range("J135").Value = 23285.8
$ws.Range("L135").Value = 23285.8
$ws.Range("N135").Value = -33425.8

# Row 136
$ws.Range("H136").Value = 1323.0741
$ws.Range("I136").Value = 1116.421
$ws.Range("J136").Value = 1813.875
$ws.Range("K136").Value = 3349.263
$ws.Range("L136").Value = 5441.625
$ws.Range("M136").Value = -799.2629999999999
$ws.Range("N136").Value = -10541.625

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1410
$ws.Range("I94").Value = 700
$ws.Range("J94").Value = 2120
$ws.Range("K94").Value = 700
$ws.Range("L94").Value = 2120
$ws.Range("M94").Value = -249
$ws.Range("N94").Value = -3022

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 16670052
$ws.Range("I31").Value = 35716830
$ws.Range("J31").Value = 4121.1875
$ws.Range("K31").Value = 35716830
$ws.Range("L31").Value = 4121.1875
$ws.Range("M31").Value = -35716535
$ws.Range("N31").Value = -4711.1875

# Row 34
$ws.Range("H34").Value = 16670052
$ws.Range("I34").Value = 35716830
$ws.Range("J34").Value = 4121.1875
$ws.Range("K34").Value = 35716830
$ws.Range("L34").Value = 4121.1875
$ws.Range("M34").Value = -35716628
$ws.Range("N34").Value = -4525.1875

# Row 58
$ws.Range("H58").Value = 1657.5
$ws.Range("J58").Value = 1579.8182
$ws.Range("L58").Value = 1579.8182
$ws.Range("N58").Value = -1985.8182

# Row 59
$ws.Range("H59").Value = 31500
$ws.Range("J59").Value = 31500
$ws.Range("L59").Value = 31500
$ws.Range("N59").Value = -33790

# Row 134
$ws.Range("H134").Value = 1900
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").Value = 5700
$ws.Range("N134").Value = -10770

# Row 136
$ws.Range("H136").Value = 1657.5
$ws.Range("J136").Value = 1579.8182
$ws.Range("L136").Value = 4739.4546
$ws.Range("N136").Value = -9839.454600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 3500
$ws.Range("J59").Value = 3500
$ws.Range("L59").Value = 10500
$ws.Range("N59").Value = -11580

# Row 121
$ws.Range("H121").Value = 1546802.5
$ws.Range("J121").Value = 1740101.5
$ws.Range("L121").Value = 5220304.5
$ws.Range("N121").Value = -5222924.5

# Row 126
$ws.Range("H126").Value = 5082.875
$ws.Range("J126").Value = 5633.2856
$ws.Range("L126").Value = 16899.8568
$ws.Range("N126").Value = -26779.8568

# Row 131
$ws.Range("H131").Value = 47337.07
$ws.Range("J131").Value = 58064.688
$ws.Range("L131").Value = 174194.064
$ws.Range("N131").Value = -184274.064

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 5908.857
$ws.Range("I132").Value = 9935.23
$ws.Range("J132").Value = 2419.3333
$ws.Range("K132").Value = 29805.69
$ws.Range("L132").Value = 7257.999899999999
$ws.Range("M132").Value = -27275.69
$ws.Range("N132").Value = -12317.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 242.73914
$ws.Range("I55").Value = 171.38095
$ws.Range("J55").Value = 992
$ws.Range("K55").Value = 171.38095
$ws.Range("L55").Value = 992
$ws.Range("M55").Value = 1.619049999999987
$ws.Range("N55").Value = -1338

# Row 93
$ws.Range("H93").Value = 2030.7391
$ws.Range("I93").Value = 1611.5
$ws.Range("J93").Value = 3540
$ws.Range("K93").Value = 1611.5
$ws.Range("L93").Value = 3540
$ws.Range("M93").Value = -363.5
$ws.Range("N93").Value = -6036

# Row 100
$ws.Range("H100").Value = 2095.4285
$ws.Range("I100").Value = 1667
$ws.Range("J100").Value = 2666.6667
$ws.Range("K100").Value = 1667
$ws.Range("L100").Value = 2666.6667
$ws.Range("M100").Value = -1126
$ws.Range("N100").Value = -3748.6667

# Row 132
$ws.Range("H132").Value = 3780.5
$ws.Range("I132").Value = 2097
$ws.Range("K132").Value = 6291
$ws.Range("M132").Value = -3761

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

Write-Host "Applied all changes"
